$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Prime the date-formatted cells (C69:D70) with the same number format
# already used by the existing date column (copy format from row 68)
# so that no new style/numFmt entries get created.
$ws.Cells.Item(68, 3).Copy()
$ws.Cells.Item(69, 3).PasteSpecial(-4122)
$ws.Cells.Item(68, 4).Copy()
$ws.Cells.Item(69, 4).PasteSpecial(-4122)
$ws.Cells.Item(68, 3).Copy()
$ws.Cells.Item(70, 3).PasteSpecial(-4122)
$ws.Cells.Item(68, 4).Copy()
$ws.Cells.Item(70, 4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 69: The Real Crash
$ws.Cells.Item(69, 1).Value = "The Real Crash"
$ws.Cells.Item(69, 2).Value = "Peter Schiff"
$ws.Cells.Item(69, 3).Value = 43952
$ws.Cells.Item(69, 4).Value = 43954
$ws.Cells.Item(69, 5).Value = "economics;national debt;politics;libertarianism;investing"
$ws.Cells.Item(69, 6).Value = "Audio"
$ws.Cells.Item(69, 7).Value = "12 Hours 55 Mins"

# Row 70: The Story of Neuroscience
$ws.Cells.Item(70, 1).Value = "The Story of Neuroscience"
$ws.Cells.Item(70, 2).Value = "Anne Rooney"
$ws.Cells.Item(70, 3).Value = 43952
$ws.Cells.Item(70, 4).Value = 43955
$ws.Cells.Item(70, 5).Value = "neuroscience;science;brain;history"
$ws.Cells.Item(70, 6).Value = "Hard Copy"
$ws.Cells.Item(70, 7).Value = "202 Pages"

# Scroll the view down and move the active selection to the row following
# the newly-added entries (best-effort; mirrors the author's on-screen state).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 51
$ws.Range("A71").Select()
